# fixed damage property type: fixed property money -> diamond
# Updates MAXHP (E) / MAXMP (F) columns for rows 11-70 (row 21 unchanged)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    11 = 550
    12 = 600
    13 = 650
    14 = 700
    15 = 750
    16 = 800
    17 = 850
    18 = 900
    19 = 950
    20 = 1000
    22 = 1100
    23 = 1150
    24 = 1200
    25 = 1250
    26 = 1300
    27 = 1350
    28 = 1400
    29 = 1450
    30 = 1500
    31 = 1550
    32 = 1600
    33 = 1650
    34 = 1700
    35 = 1750
    36 = 1800
    37 = 1850
    38 = 1900
    39 = 1950
    40 = 2000
    41 = 2050
    42 = 2100
    43 = 2150
    44 = 2200
    45 = 2250
    46 = 2300
    47 = 2350
    48 = 2400
    49 = 2450
    50 = 2500
    51 = 2550
    52 = 2600
    53 = 2650
    54 = 2700
    55 = 2750
    56 = 2800
    57 = 2850
    58 = 2900
    59 = 2950
    60 = 3000
    61 = 3050
    62 = 3100
    63 = 3150
    64 = 3200
    65 = 3250
    66 = 3300
    67 = 3350
    68 = 3400
    69 = 3450
    70 = 3500
}

foreach ($row in $newValues.Keys) {
    $val = $newValues[$row]
    $ws.Range("E$row").Value = $val
    $ws.Range("F$row").Value = $val
}

# Restore the active selection to match the saved workbook state
$ws.Range("G30").Select()
